# Apply updated IPC PO predictions and derived DELTA / DELTA^2 values
# following the refactor of weight handling in DenseLayer/NeuralNetwork.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 29.36188838516495
$ws.Range("D2").Value = -0.1781116148350463
$ws.Range("E2").Value = 0.0317237473391479
$ws.Range("C3").Value = 29.78406929991331
$ws.Range("D3").Value = 0.2340692999133118
$ws.Range("E3").Value = 0.05478843716190792
$ws.Range("C4").Value = 29.73033275421328
$ws.Range("D4").Value = -0.01966724578672441
$ws.Range("E4").Value = 0.000386800556835429
$ws.Range("C5").Value = 30.29087302102019
$ws.Range("D5").Value = 0.4508730210201932
$ws.Range("E5").Value = 0.2032864810838756
$ws.Range("C6").Value = 30.16897008820004
$ws.Range("D6").Value = 0.3589700882000386
$ws.Range("E6").Value = 0.1288595242223435
$ws.Range("C7").Value = 29.90822042436572
$ws.Range("D7").Value = -0.01177957563428578
$ws.Range("E7").Value = 0.0001387584021238592
$ws.Range("C8").Value = 29.83849563355738
$ws.Range("D8").Value = -0.1415043664426214
$ws.Range("E8").Value = 0.02002348572232769
$ws.Range("C9").Value = 29.84181641502009
$ws.Range("D9").Value = -0.1981835849799047
$ws.Range("E9").Value = 0.03927673335548712
$ws.Range("C10").Value = 29.92971577386612
$ws.Range("D10").Value = -0.2802842261338796
$ws.Range("E10").Value = 0.07855924741946775
$ws.Range("C11").Value = 30.09036226924875
$ws.Range("D11").Value = -0.1296377307512451
$ws.Range("E11").Value = 0.01680594123433231
$ws.Range("C12").Value = 30.19173587840176
$ws.Range("D12").Value = -0.188264121598241
$ws.Range("E12").Value = 0.03544337948115728
$ws.Range("C13").Value = 30.48602904936035
$ws.Range("D13").Value = 0.04602904936034946
$ws.Range("E13").Value = 0.002118673385017487
$ws.Range("C14").Value = 30.3830970472716
$ws.Range("D14").Value = -0.09690295272840288
$ws.Range("E14").Value = 0.009390182247483083
$ws.Range("C15").Value = 30.40061040351136
$ws.Range("D15").Value = -0.2893895964886433
$ws.Range("E15").Value = 0.08374633855585976
$ws.Range("C16").Value = 30.75093117328817
$ws.Range("D16").Value = 0.0009311732881656098
$ws.Range("E16").Value = 0.000000867083692593153831567135
$ws.Range("C17").Value = 30.80748496167179
$ws.Range("D17").Value = -0.1325150383282079
$ws.Range("E17").Value = 0.0175602353831264
$ws.Range("C18").Value = 30.92476076942891
$ws.Range("D18").Value = -0.02523923057109201
$ws.Range("E18").Value = 0.0006370187598207454
$ws.Range("C19").Value = 31.14990447558944
$ws.Range("D19").Value = 0.1299044755894379
$ws.Range("E19").Value = 0.01687517277816685
$ws.Range("C20").Value = 31.28542580370341
$ws.Range("D20").Value = 0.1654258037034069
$ws.Range("E20").Value = 0.0273656965309181
$ws.Range("C21").Value = 31.1840818168263
$ws.Range("D21").Value = -0.09591818317370127
$ws.Range("E21").Value = 0.009200297863343709
$ws.Range("C22").Value = 31.07793052009108
$ws.Range("D22").Value = -0.3020694799089227
$ws.Range("E22").Value = 0.09124597069244708
$ws.Range("C23").Value = 31.27197894129191
$ws.Range("D23").Value = -0.3080210587080927
$ws.Range("E23").Value = 0.09487697260765429
$ws.Range("C24").Value = 31.73957506555968
$ws.Range("D24").Value = 0.08957506555968564
$ws.Range("E24").Value = 0.008023692370021981
$ws.Range("C25").Value = 32.50872068394579
$ws.Range("D25").Value = 0.6287206839457902
$ws.Range("E25").Value = 0.3952896984212622
$ws.Range("C26").Value = 32.53122384196269
$ws.Range("D26").Value = 0.251223841962684
$ws.Range("E26").Value = 0.06311341877049162
$ws.Range("C27").Value = 32.88905978555466
$ws.Range("D27").Value = 0.4390597855546616
$ws.Range("E27").Value = 0.1927734952913055
$ws.Range("C28").Value = 33.12466019119131
$ws.Range("D28").Value = 0.2746601911913089
$ws.Range("E28").Value = 0.07543822062524637
$ws.Range("C29").Value = 33.1876042672486
$ws.Range("D29").Value = 0.2876042672486037
$ws.Range("E29").Value = 0.08271621453960623
$ws.Range("C30").Value = 33.31342044662783
$ws.Range("D30").Value = 0.2134204466278291
$ws.Range("E30").Value = 0.04554828703882205
$ws.Range("C31").Value = 33.71395351623086
$ws.Range("D31").Value = 0.3139535162308604
$ws.Range("E31").Value = 0.09856681035372114
$ws.Range("C32").Value = 33.72374872713567
$ws.Range("D32").Value = 0.02374872713566845
$ws.Range("E32").Value = 0.0005640020405644348
$ws.Range("C33").Value = 33.79202561111252
$ws.Range("D33").Value = -0.3079743888874802
$ws.Range("E33").Value = 0.09484822421061688
$ws.Range("C34").Value = 34.19723309020993
$ws.Range("D34").Value = -0.2027669097900642
$ws.Range("E34").Value = 0.04111441970581203
$ws.Range("C35").Value = 34.41133269554031
$ws.Range("D35").Value = -0.488667304459689
$ws.Range("E35").Value = 0.2387957344478984
$ws.Range("C36").Value = 35.51423373159818
$ws.Range("D36").Value = 0.2142337315981848
$ws.Range("E36").Value = 0.04589609175448309
$ws.Range("C37").Value = 35.77938616129119
$ws.Range("D37").Value = 0.07938616129118259
$ws.Range("E37").Value = 0.006302162604549657
$ws.Range("C38").Value = 36.05415544559393
$ws.Range("D38").Value = -0.2458445544060694
$ws.Range("E38").Value = 0.06043954493111884
$ws.Range("C39").Value = 36.66889473158334
$ws.Range("D39").Value = -0.131105268416654
$ws.Range("E39").Value = 0.0171885914066029
$ws.Range("C40").Value = 36.75928092501331
$ws.Range("D40").Value = -0.5407190749866899
$ws.Range("E40").Value = 0.2923771180544616
$ws.Range("C41").Value = 37.85027403718746
$ws.Range("D41").Value = -0.04972596281253772
$ws.Range("E41").Value = 0.002472671377633884
$ws.Range("C42").Value = 38.61013807581099
$ws.Range("D42").Value = 0.110138075810994
$ws.Range("E42").Value = 0.01213039574334826
$ws.Range("C43").Value = 39.00815598828223
$ws.Range("D43").Value = 0.1081559882822347
$ws.Range("E43").Value = 0.01169771780130689
$ws.Range("C44").Value = 39.39532542033854
$ws.Range("D44").Value = -0.004674579661461564
$ws.Range("E44").Value = 0.000021851695011350108515833443
$ws.Range("C45").Value = 39.65916224730229
$ws.Range("D45").Value = -0.2408377526977077
$ws.Range("E45").Value = 0.0580028231244822
$ws.Range("C46").Value = 39.93191622206452
$ws.Range("D46").Value = -0.1680837779354789
$ws.Range("E46").Value = 0.02825215640506339
$ws.Range("C47").Value = 40.14516499738765
$ws.Range("D47").Value = -0.4548350026123487
$ws.Range("E47").Value = 0.2068748796013752
$ws.Range("C48").Value = 40.46227174547302
$ws.Range("D48").Value = -0.4377282545269807
$ws.Range("E48").Value = 0.1916060248112372
$ws.Range("C49").Value = 41.57788098621841
$ws.Range("D49").Value = 0.3778809862184076
$ws.Range("E49").Value = 0.1427940397453964
$ws.Range("C50").Value = 41.69234694904898
$ws.Range("D50").Value = 0.1923469490489822
$ws.Range("E50").Value = 0.03699734880845176
$ws.Range("C51").Value = 42.12059394226941
$ws.Range("D51").Value = 0.3205939422694115
$ws.Range("E51").Value = 0.1027804758198428
$ws.Range("C52").Value = -0.3595455662107803
$ws.Range("E52").Value = 3.51493607336627
$ws.Range("E53").Value = 0.07029872146732541
